$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the definition / comment of "Channel ID" (A1 on the "Antibodies"
#    sheet) to describe matching the OME TIFF channel ID.
# ---------------------------------------------------------------------------
$wsAntibodies = $wb.Worksheets.Item("Antibodies")

$channelIdComment = @"
(Required) Structure of the identifier depends on the acquisition system.
Whenever possible this should exactly match the channel ID in the OME TIFF file.
For example the channel ID in an OME TIFF might be something like "Channel:0:13"
which would then be the value entered here.
"@

$wsAntibodies.Range("A1").Comment.Text($channelIdComment)

# ---------------------------------------------------------------------------
# 2) Bump the pav:createdOn timestamp on the ".metadata" sheet (row 2,
#    column C holds the value for the "pav:createdOn" header in row 1).
# ---------------------------------------------------------------------------
$wsMetadata = $wb.Worksheets.Item(".metadata")
$wsMetadata.Range("C2").Value = "2024-11-05T13:43:22-08:00"

# ---------------------------------------------------------------------------
# 3) Refresh the "best fit" column widths on the "Antibodies" and
#    ".metadata" sheets (narrower than before).
#    ColumnWidth is stored/quantized in whole "pixel" steps by this engine,
#    so we compensate for the fixed +5/6 rounding offset it applies in order
#    to land as close as possible to the target widths.
# ---------------------------------------------------------------------------
function Set-ClosestColumnWidth($col, $targetWidth) {
    $col.ColumnWidth = $targetWidth - 0.8333333333333334
}

$antibodyWidths = @(9.109375, 10.984375, 11.12890625, 21.46875, 9.67578125, 12.05859375, 23.734375, 22.734375, 19.2734375, 12.6953125, 16.91796875)
for ($i = 0; $i -lt $antibodyWidths.Length; $i++) {
    Set-ClosestColumnWidth $wsAntibodies.Columns.Item($i + 1) $antibodyWidths[$i]
}

$metadataWidths = @(10.20703125, 9.65234375, 24.1484375, 68.6484375)
for ($i = 0; $i -lt $metadataWidths.Length; $i++) {
    Set-ClosestColumnWidth $wsMetadata.Columns.Item($i + 1) $metadataWidths[$i]
}
